$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.998.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.85%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.718.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.39%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.81%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.72%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.22%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.578'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.71%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.744.52'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.07%  '

# Row 10
$ws.Range("E10").Value = '  +12.00%  '

# Row 11
$ws.Range("E11").Value = '  +1.95%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.344'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.75%  '

# Row 13
$ws.Range("E13").Value = '  +2.67%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.204.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.20%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.955.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.77%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.27%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.779.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.46%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000140'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.15%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '347.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.10%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.84%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.25%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.50%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.06%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.173'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.80%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.423'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.11%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.995'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.38%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0839'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.30%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.07%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.10%  '

# Row 31
$ws.Range("E31").Value = '  -0.14%  '

# Row 32
$ws.Range("E32").Value = '  +2.67%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.94%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '151.61'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.38%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.63%  '

# Row 36
$ws.Range("E36").Value = '  +7.81%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.922'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.13%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.914'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.49%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.60%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.62'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.92%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.42%  '

# Row 42
$ws.Range("E42").Value = '  +4.72%  '

# Row 43
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '283.63'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.32%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.35'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.38%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.996'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.13%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0987'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.16%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.28%  '

# Row 48
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0549'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.55%  '

# Row 49
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.109.79'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.85%  '

# Row 50
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.73'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.95%  '

# Row 51
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.24%  '
